$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row (columns A..T) for sheet rows 2..7
$rowsData = @(
    ,("ECs", "Fn1", "Col13a1", "ECs", 3, 1, 21.84976866666667, 65.549306, 0.05020018890879543, 0.05020018890879543, 2, 0.6666666666666666, 0.3326186666666667, 0.9978560000000001, 0.3213922220876632, 0.3213922220876632, 7.267640920881778, 65.40876828793601, 0.01613395026261823, 0.01613395026261823)
    ,("ECs", "Fn1", "Col13a1", "sCs", 3, 1, 21.84976866666667, 65.549306, 0.05020018890879543, 0.05020018890879543, 3, 1, 0.702312, 2.106936, 0.6786077779123368, 0.6786077779123368, 15.345354731824, 138.108192586416, 0.0340662386461772, 0.0340662386461772)
    ,("FAPs", "Fn1", "Col13a1", "ECs", 3, 1, 385.0524703333334, 1155.157411, 0.8846641374295412, 0.8846641374295412, 2, 0.6666666666666666, 0.3326186666666667, 0.9978560000000001, 0.3213922220876632, 0.3213922220876632, 128.0756392789796, 1152.680753510816, 0.2843241729297461, 0.2843241729297461)
    ,("FAPs", "Fn1", "Col13a1", "sCs", 3, 1, 385.0524703333334, 1155.157411, 0.8846641374295412, 0.8846641374295412, 3, 1, 0.702312, 2.106936, 0.6786077779123368, 0.6786077779123368, 270.4269705447441, 2433.842734902697, 0.6003399644997951, 0.6003399644997951)
    ,("sCs", "Fn1", "Col13a1", "ECs", 3, 1, 28.350479, 85.05143699999999, 0.06513567366166337, 0.06513567366166337, 2, 0.6666666666666666, 0.3326186666666667, 0.9978560000000001, 0.3213922220876632, 0.3213922220876632, 9.429898524341333, 84.869086719072, 0.02093409889529887, 0.02093409889529887)
    ,("sCs", "Fn1", "Col13a1", "sCs", 3, 1, 28.350479, 85.05143699999999, 0.06513567366166337, 0.06513567366166337, 3, 1, 0.702312, 2.106936, 0.6786077779123368, 0.6786077779123368, 19.910881607448, 179.197934467032, 0.0442015747663645, 0.0442015747663645)
)

$numRows = $rowsData.Count
$numCols = $rowsData[0].Count
$arr = New-Object 'object[,]' $numRows,$numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $arr[$r,$c] = $rowsData[$r][$c]
    }
}

$ws.Range("A2:T7").Value = $arr